$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 389, pushing existing data (389..461) down to (392..464)
$ws.Rows.Item(389).Resize(3).Insert()

# Common (static) values shared by every data row in this block
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100101
$producto  = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$variedad  = "Hayward"
$unidad    = "$/caja 15 kilos"
$origen    = "Región de O'Higgins"
$kgUnidad  = 15

# New row 389 - Especial
$ws.Cells.Item(389, 1).Value = $mercadoId
$ws.Cells.Item(389, 2).Value = $mercado
$ws.Cells.Item(389, 3).Value = $region
$ws.Cells.Item(389, 4).Value = 45015
$ws.Cells.Item(389, 5).Value = $codreg
$ws.Cells.Item(389, 6).Value = $tipo
$ws.Cells.Item(389, 7).Value = $productoId
$ws.Cells.Item(389, 8).Value = $producto
$ws.Cells.Item(389, 9).Value = $categoriaId
$ws.Cells.Item(389, 10).Value = $categoria
$ws.Cells.Item(389, 11).Value = $variedad
$ws.Cells.Item(389, 12).Value = "Especial"
$ws.Cells.Item(389, 13).Value = 250
$ws.Cells.Item(389, 14).Value = 22000
$ws.Cells.Item(389, 15).Value = 22000
$ws.Cells.Item(389, 16).Value = 22000
$ws.Cells.Item(389, 17).Value = $unidad
$ws.Cells.Item(389, 18).Value = $origen
$ws.Cells.Item(389, 19).Value = 1467
$ws.Cells.Item(389, 20).Value = $kgUnidad

# New row 390 - Primera
$ws.Cells.Item(390, 1).Value = $mercadoId
$ws.Cells.Item(390, 2).Value = $mercado
$ws.Cells.Item(390, 3).Value = $region
$ws.Cells.Item(390, 4).Value = 45015
$ws.Cells.Item(390, 5).Value = $codreg
$ws.Cells.Item(390, 6).Value = $tipo
$ws.Cells.Item(390, 7).Value = $productoId
$ws.Cells.Item(390, 8).Value = $producto
$ws.Cells.Item(390, 9).Value = $categoriaId
$ws.Cells.Item(390, 10).Value = $categoria
$ws.Cells.Item(390, 11).Value = $variedad
$ws.Cells.Item(390, 12).Value = "Primera"
$ws.Cells.Item(390, 13).Value = 250
$ws.Cells.Item(390, 14).Value = 19000
$ws.Cells.Item(390, 15).Value = 19000
$ws.Cells.Item(390, 16).Value = 19000
$ws.Cells.Item(390, 17).Value = $unidad
$ws.Cells.Item(390, 18).Value = $origen
$ws.Cells.Item(390, 19).Value = 1267
$ws.Cells.Item(390, 20).Value = $kgUnidad

# New row 391 - Segunda
$ws.Cells.Item(391, 1).Value = $mercadoId
$ws.Cells.Item(391, 2).Value = $mercado
$ws.Cells.Item(391, 3).Value = $region
$ws.Cells.Item(391, 4).Value = 45015
$ws.Cells.Item(391, 5).Value = $codreg
$ws.Cells.Item(391, 6).Value = $tipo
$ws.Cells.Item(391, 7).Value = $productoId
$ws.Cells.Item(391, 8).Value = $producto
$ws.Cells.Item(391, 9).Value = $categoriaId
$ws.Cells.Item(391, 10).Value = $categoria
$ws.Cells.Item(391, 11).Value = $variedad
$ws.Cells.Item(391, 12).Value = "Segunda"
$ws.Cells.Item(391, 13).Value = 250
$ws.Cells.Item(391, 14).Value = 16000
$ws.Cells.Item(391, 15).Value = 16000
$ws.Cells.Item(391, 16).Value = 16000
$ws.Cells.Item(391, 17).Value = $unidad
$ws.Cells.Item(391, 18).Value = $origen
$ws.Cells.Item(391, 19).Value = 1067
$ws.Cells.Item(391, 20).Value = $kgUnidad
